# Auto-generated edit script applying the Kujata_Profits market-price refresh
# across all profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 1619.8334
$ws.Range("I70").Value = 1489.3334
$ws.Range("J70").Value = 1663.3334
$ws.Range("K70").Value = 4468.0002
$ws.Range("L70").Value = 4990.0002
$ws.Range("M70").Value = -4198.0002
$ws.Range("N70").Value = -5530.0002
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 1619.8334
$ws.Range("I73").Value = 1489.3334
$ws.Range("J73").Value = 1663.3334
$ws.Range("K73").Value = 4468.0002
$ws.Range("L73").Value = 4990.0002
$ws.Range("M73").Value = -3532.0002
$ws.Range("N73").Value = -6862.0002
$ws.Range("H86").Value = 6543.4287
$ws.Range("I86").Value = 7200
$ws.Range("J86").Value = 4902
$ws.Range("K86").Value = 7200
$ws.Range("L86").Value = 4902
$ws.Range("M86").Value = -6077
$ws.Range("N86").Value = -7148
$ws.Range("H89").Value = 6543.4287
$ws.Range("I89").Value = 7200
$ws.Range("J89").Value = 4902
$ws.Range("K89").Value = 36000
$ws.Range("L89").Value = 24510
$ws.Range("M89").Value = -30384
$ws.Range("N89").Value = -35742
$ws.Range("H112").Value = 1932.65
$ws.Range("J112").Value = 1987.1052
$ws.Range("L112").Value = 5961.3156
$ws.Range("N112").Value = -8177.3156
$ws.Range("H129").Value = 910.25
$ws.Range("I129").Value = 821.25
$ws.Range("J129").Value = 999.25
$ws.Range("K129").Value = 2463.75
$ws.Range("L129").Value = 2997.75
$ws.Range("M129").Value = 2536.25
$ws.Range("N129").Value = -12997.75
$ws.Range("H137").Value = 2386.7334
$ws.Range("I137").Value = 6499
$ws.Range("J137").Value = 1754.0769
$ws.Range("K137").Value = 19497
$ws.Range("L137").Value = 5262.2307
$ws.Range("M137").Value = -16947
$ws.Range("N137").Value = -10362.2307
$ws.Range("H138").Value = 3566.5466
$ws.Range("I138").Value = 2000.5
$ws.Range("J138").Value = 3925.9673
$ws.Range("K138").Value = 6001.5
$ws.Range("L138").Value = 11777.9019
$ws.Range("M138").Value = -861.5
$ws.Range("N138").Value = -22057.9019

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15595.978
$ws.Range("I32").Value = 14195.81
$ws.Range("K32").Value = 14195.81
$ws.Range("M32").Value = -13908.81
$ws.Range("H40").Value = 6500
$ws.Range("J40").Value = 6500
$ws.Range("L40").Value = 6500
$ws.Range("N40").Value = -6852
$ws.Range("H74").Value = 2609.4
$ws.Range("I74").Value = 1107.5
$ws.Range("K74").Value = 1107.5
$ws.Range("M74").Value = -233.5
$ws.Range("H77").Value = 2609.4
$ws.Range("I77").Value = 1107.5
$ws.Range("K77").Value = 5537.5
$ws.Range("M77").Value = -1169.5
$ws.Range("H88").Value = 1743.3636
$ws.Range("I88").Value = 1622.5
$ws.Range("J88").Value = 1812.4286
$ws.Range("K88").Value = 1622.5
$ws.Range("L88").Value = 1812.4286
$ws.Range("M88").Value = -1216.5
$ws.Range("N88").Value = -2624.4286
$ws.Range("H91").Value = 1743.3636
$ws.Range("I91").Value = 1622.5
$ws.Range("J91").Value = 1812.4286
$ws.Range("K91").Value = 1622.5
$ws.Range("L91").Value = 1812.4286
$ws.Range("M91").Value = -218.5
$ws.Range("N91").Value = -4620.4286
$ws.Range("H122").Value = 4367.684
$ws.Range("J122").Value = 4473.778
$ws.Range("L122").Value = 13421.334
$ws.Range("N122").Value = -18321.334
$ws.Range("H132").Value = 2860.4866
$ws.Range("I132").Value = 2057.1
$ws.Range("J132").Value = 3805.647
$ws.Range("K132").Value = 6171.299999999999
$ws.Range("L132").Value = 11416.941
$ws.Range("M132").Value = -3641.299999999999
$ws.Range("N132").Value = -16476.941

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H42").Value = 132342
$ws.Range("J42").Value = 132342
$ws.Range("L42").Value = 132342
$ws.Range("N42").Value = -132998
$ws.Range("H86").Value = 3874.3044
$ws.Range("I86").Value = 3872.889
$ws.Range("J86").Value = 3879.4
$ws.Range("K86").Value = 3872.889
$ws.Range("L86").Value = 3879.4
$ws.Range("M86").Value = -2749.889
$ws.Range("N86").Value = -6125.4
$ws.Range("H89").Value = 3874.3044
$ws.Range("I89").Value = 3872.889
$ws.Range("J89").Value = 3879.4
$ws.Range("K89").Value = 19364.445
$ws.Range("L89").Value = 19397
$ws.Range("M89").Value = -13748.445
$ws.Range("N89").Value = -30629
$ws.Range("H134").Value = 3185.4082
$ws.Range("I134").Value = 892.4516
$ws.Range("J134").Value = 7134.3887
$ws.Range("K134").Value = 2677.3548
$ws.Range("L134").Value = 21403.1661
$ws.Range("M134").Value = -142.3548000000001
$ws.Range("N134").Value = -26473.1661

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 33000
$ws.Range("J74").Value = 33000
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34748
$ws.Range("H77").Value = 33000
$ws.Range("J77").Value = 33000
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -107736
$ws.Range("H132").Value = 2052.2104
$ws.Range("I132").Value = 1825.6552
$ws.Range("J132").Value = 2782.2222
$ws.Range("K132").Value = 5476.9656
$ws.Range("L132").Value = 8346.6666
$ws.Range("M132").Value = -2946.9656
$ws.Range("N132").Value = -13406.6666
$ws.Range("H141").Value = 494602.78
$ws.Range("J141").Value = 494602.78
$ws.Range("L141").Value = 494602.78
$ws.Range("N141").Value = -504962.78

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2180.4
$ws.Range("I68").Value = 1851
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 5553
$ws.Range("L68").Value = 7200
$ws.Range("M68").Value = -4742
$ws.Range("N68").Value = -8822
$ws.Range("H71").Value = 2180.4
$ws.Range("I71").Value = 1851
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 16659
$ws.Range("L71").Value = 21600
$ws.Range("M71").Value = -12603
$ws.Range("N71").Value = -29712
$ws.Range("H107").Value = 6806.579
$ws.Range("I107").Value = 372.22223
$ws.Range("K107").Value = 1116.66669
$ws.Range("M107").Value = 803.33331
$ws.Range("H113").Value = 721.7143
$ws.Range("J113").Value = 721.7143
$ws.Range("L113").Value = 2165.1429
$ws.Range("N113").Value = -6505.1429

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 53400
$ws.Range("J104").Value = 53400
$ws.Range("L104").Value = 53400
$ws.Range("N104").Value = -60388
$ws.Range("I122").Value = 1344.8572
$ws.Range("J122").Value = 1449.6666
$ws.Range("K122").Value = 4034.5716
$ws.Range("L122").Value = 4348.9998
$ws.Range("M122").Value = -1584.5716
$ws.Range("N122").Value = -9248.9998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 25759810
$ws.Range("I122").Value = 56668040
$ws.Range("J122").Value = 2951.5
$ws.Range("K122").Value = 170004120
$ws.Range("L122").Value = 8854.5
$ws.Range("M122").Value = -170001670
$ws.Range("N122").Value = -13754.5
$ws.Range("H136").Value = 1964.9048
$ws.Range("I136").Value = 1741.125
$ws.Range("K136").Value = 5223.375
$ws.Range("M136").Value = -2673.375
$ws.Range("H139").Value = 50715
$ws.Range("J139").Value = 50715
$ws.Range("L139").Value = 50715
$ws.Range("N139").Value = -60995
$ws.Range("H140").Value = 51158.9
$ws.Range("J140").Value = 51158.9
$ws.Range("L140").Value = 51158.9
$ws.Range("N140").Value = -61518.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4851
$ws.Range("H62").Value = 62503310
$ws.Range("I62").Value = 100002800
$ws.Range("K62").Value = 100002800
$ws.Range("M62").Value = -100002176
$ws.Range("H65").Value = 62503310
$ws.Range("I65").Value = 100002800
$ws.Range("K65").Value = 500014000
$ws.Range("M65").Value = -500010880
$ws.Range("H81").Value = 1840.9546
$ws.Range("I81").Value = 1175.25
$ws.Range("J81").Value = 1988.8889
$ws.Range("K81").Value = 2350.5
$ws.Range("L81").Value = 3977.7778
$ws.Range("M81").Value = -1289.5
$ws.Range("N81").Value = -6099.7778
$ws.Range("H84").Value = 1840.9546
$ws.Range("I84").Value = 1175.25
$ws.Range("J84").Value = 1988.8889
$ws.Range("K84").Value = 11752.5
$ws.Range("L84").Value = 19888.889
$ws.Range("M84").Value = -6448.5
$ws.Range("N84").Value = -30496.889
$ws.Range("H122").Value = 7814755
$ws.Range("I122").Value = 10001979
$ws.Range("J122").Value = 3240.7144
$ws.Range("K122").Value = 30005937
$ws.Range("L122").Value = 9722.143199999999
$ws.Range("M122").Value = -30003487
$ws.Range("N122").Value = -14622.1432
$ws.Range("H126").Value = 50506076
$ws.Range("I126").Value = 85471040
$ws.Range("J126").Value = 1126.6666
$ws.Range("K126").Value = 256413120
$ws.Range("L126").Value = 3379.9998
$ws.Range("M126").Value = -256410650
$ws.Range("N126").Value = -8319.9998
$ws.Range("H132").Value = 4948.5625
$ws.Range("I132").Value = 5671.4
$ws.Range("J132").Value = 3743.8333
$ws.Range("K132").Value = 17014.2
$ws.Range("L132").Value = 11231.4999
$ws.Range("M132").Value = -14484.2
$ws.Range("N132").Value = -16291.4999
$ws.Range("H136").Value = 1144
$ws.Range("J136").Value = 2071.4285
$ws.Range("L136").Value = 6214.2855
$ws.Range("N136").Value = -11314.2855
